# Update "想去人数" (F column) values across sheets, matching the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 432
$ws1.Range("F6").Value = 1262
$ws1.Range("F8").Value = 7555
$ws1.Range("F10").Value = 108
$ws1.Range("F12").Value = 8201
$ws1.Range("F13").Value = 1
$ws1.Range("F16").Value = 5614
$ws1.Range("F18").Value = 2566
$ws1.Range("F25").Value = 476
$ws1.Range("F26").Value = 2526
$ws1.Range("F29").Value = 2794
$ws1.Range("F30").Value = 0
$ws1.Range("F31").Value = 322
$ws1.Range("F32").Value = 117
$ws1.Range("F33").Value = 274
$ws1.Range("F34").Value = 633
$ws1.Range("F40").Value = 2607
$ws1.Range("F42").Value = 2268

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 96
$ws2.Range("F6").Value = 33

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1305

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1305
$ws4.Range("F6").Value = 1262
$ws4.Range("F7").Value = 7555
$ws4.Range("F9").Value = 108
$ws4.Range("F11").Value = 8201
$ws4.Range("F12").Value = 1
$ws4.Range("F15").Value = 5614
$ws4.Range("F17").Value = 2566
$ws4.Range("F22").Value = 96
$ws4.Range("F25").Value = 476
$ws4.Range("F26").Value = 2526
$ws4.Range("F29").Value = 2795
$ws4.Range("F30").Value = 322
$ws4.Range("F31").Value = 117
$ws4.Range("F32").Value = 274
$ws4.Range("F34").Value = 633
$ws4.Range("F38").Value = 33
$ws4.Range("F42").Value = 2607
$ws4.Range("F45").Value = 2268
